# Apply cryptos list price/volume updates (and two name/link swaps) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.845.20"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "2.739.98"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'572.45"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("D6").Value = "'156.44"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.600"
$ws.Range("E8").Value = "  -1.61%  "

$ws.Range("D9").Value = "'0.109"
$ws.Range("E9").Value = "  -3.18%  "

$ws.Range("E10").Value = "  -1.00%  "

$ws.Range("D11").Value = "'0.380"
$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "'5.57"
$ws.Range("E12").Value = "  -16.84%  "

$ws.Range("D13").Value = "3.224.41"
$ws.Range("E13").Value = "  -0.45%  "

$ws.Range("D14").Value = "'26.42"
$ws.Range("E14").Value = "  +0.52%  "

$ws.Range("D15").Value = "63.497.46"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("D16").Value = "'0.0000148"
$ws.Range("E16").Value = "  -2.73%  "

$ws.Range("D17").Value = "2.745.94"
$ws.Range("E17").Value = "  -0.34%  "

$ws.Range("D18").Value = "'12.04"
$ws.Range("E18").Value = "  +0.80%  "

$ws.Range("D19").Value = "'4.79"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").Value = "'353.86"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").Value = "'6.74"
$ws.Range("E21").Value = "  -3.25%  "

$ws.Range("D22").Value = "'0.537"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("D23").Value = "'0.995"
$ws.Range("E23").Value = "  -0.59%  "

$ws.Range("D24").Value = "'65.16"
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("D25").Value = "'0.169"
$ws.Range("E25").Value = "  -0.54%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'8.36"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").Value = "0.0₃0882"
$ws.Range("E28").Value = "  -3.75%  "

$ws.Range("D29").Value = "'1.93"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").Value = "'6.89"
$ws.Range("E30").Value = "  -2.96%  "

$ws.Range("D31").Value = "'169.07"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "'1.18"
$ws.Range("E32").Value = "  -3.71%  "

$ws.Range("D33").Value = "'20.05"
$ws.Range("E33").Value = "  -1.91%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.85"
$ws.Range("E35").Value = "  +1.37%  "

$ws.Range("D36").Value = "'1.42"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "'1.77"
$ws.Range("E37").Value = "  -2.25%  "

$ws.Range("D38").Value = "'0.967"
$ws.Range("E38").Value = "  -6.16%  "

$ws.Range("D39").Value = "'6.14"
$ws.Range("E39").Value = "  +9.65%  "

$ws.Range("D40").Value = "'4.09"
$ws.Range("E40").Value = "  -2.94%  "

$ws.Range("D41").Value = "'322.00"
$ws.Range("E41").Value = "  -6.88%  "

$ws.Range("D42").Value = "'38.81"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").Value = "'21.16"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("D44").Value = "'0.0583"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").Value = "'21.21"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("D46").Value = "'0.0252"
$ws.Range("E46").Value = "  -1.03%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.623"
$ws.Range("E47").Value = "  -3.76%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'134.27"
$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'11.05"
$ws.Range("E51").Value = "  +0.63%  "
